# fix(module3): use uncon_planned_qty for future production; keep produced for today
# The rows are re-derived (sorted by material, then location DC_001/DC_002/PLANT_001)
# with new quantity/layer/horizon values, and a new MAT_B / PLANT_001 row is appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateFormat = "YYYY-MM-DD HH:MM:SS"

# Row 2: MAT_A / DC_001
$ws.Range("A2").Value = "MAT_A"
$ws.Range("B2").Value = "DC_001"
$ws.Range("C2").Value = 45298
$ws.Range("C2").NumberFormat = $dateFormat
$ws.Range("D2").Value = "Distribution Demand - Forecast"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = -239
$ws.Range("G2").Value = 45297
$ws.Range("G2").NumberFormat = $dateFormat
$ws.Range("H2").Value = 4

# Row 3: MAT_A / DC_002
$ws.Range("A3").Value = "MAT_A"
$ws.Range("B3").Value = "DC_002"
$ws.Range("C3").Value = 45298
$ws.Range("C3").NumberFormat = $dateFormat
$ws.Range("D3").Value = "Distribution Demand - Forecast"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = -562
$ws.Range("G3").Value = 45297
$ws.Range("G3").NumberFormat = $dateFormat
$ws.Range("H3").Value = 4

# Row 4: MAT_A / PLANT_001 (was MAT_B / DC_001)
$ws.Range("A4").Value = "MAT_A"
$ws.Range("B4").Value = "PLANT_001"
$ws.Range("C4").Value = 45298
$ws.Range("C4").NumberFormat = $dateFormat
$ws.Range("D4").Value = "Distribution Demand - Forecast"
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = -863
$ws.Range("G4").Value = 45297
$ws.Range("G4").NumberFormat = $dateFormat
$ws.Range("H4").Value = 1

# Row 5: MAT_B / DC_001 (was MAT_B / DC_002)
$ws.Range("A5").Value = "MAT_B"
$ws.Range("B5").Value = "DC_001"
$ws.Range("C5").Value = 45298
$ws.Range("C5").NumberFormat = $dateFormat
$ws.Range("D5").Value = "Distribution Demand - Forecast"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = -113
$ws.Range("G5").Value = 45297
$ws.Range("G5").NumberFormat = $dateFormat
$ws.Range("H5").Value = 4

# Row 6: MAT_B / DC_002 (was MAT_B / PLANT_001)
$ws.Range("A6").Value = "MAT_B"
$ws.Range("B6").Value = "DC_002"
$ws.Range("C6").Value = 45298
$ws.Range("C6").NumberFormat = $dateFormat
$ws.Range("D6").Value = "Distribution Demand - Forecast"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = -32
$ws.Range("G6").Value = 45297
$ws.Range("G6").NumberFormat = $dateFormat
$ws.Range("H6").Value = 1

# Row 7 (new): MAT_B / PLANT_001
$ws.Range("A7").Value = "MAT_B"
$ws.Range("B7").Value = "PLANT_001"
$ws.Range("C7").Value = 45298
$ws.Range("C7").NumberFormat = $dateFormat
$ws.Range("D7").Value = "Distribution Demand - Forecast"
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = -113
$ws.Range("G7").Value = 45297
$ws.Range("G7").NumberFormat = $dateFormat
$ws.Range("H7").Value = 1
